$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must stay text (matches original
# inline-string typing of the Price column), so force Text format before writing.
$textForceCells = @('D5', 'D6', 'D7', 'D10', 'D14', 'D15', 'D17', 'D18', 'D23', 'D24', 'D25', 'D27', 'D30', 'D32', 'D33', 'D34', 'D36', 'D37', 'D38', 'D39', 'D41', 'D43', 'D44', 'D46', 'D50', 'D51')
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values scraped for this run.
$ws.Range('D2').Value = '51.447.23'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '3.090.08'
$ws.Range('E3').Value = '  +2.28%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '384.32'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('D6').Value = '102.35'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').Value = '0.539'
$ws.Range('E7').Value = '  -0.93%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('D10').Value = '36.98'
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '3.579.64'
$ws.Range('E13').Value = '  +2.23%  '
$ws.Range('D14').Value = '18.59'
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').Value = '7.80'
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').Value = '3.093.98'
$ws.Range('E16').Value = '  +2.29%  '
$ws.Range('D17').Value = '11.18'
$ws.Range('E17').Value = '  +5.49%  '
$ws.Range('D18').Value = '0.990'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').Value = '51.435.81'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('E20').Value = '  +8.57%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  -0.93%  '
$ws.Range('D23').Value = '69.86'
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('D24').Value = '265.28'
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').Value = '3.12'
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('E26').Value = '  -1.46%  '
$ws.Range('D27').Value = '26.92'
$ws.Range('E27').Value = '  +2.77%  '
$ws.Range('E28').Value = '  -2.76%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '0.166'
$ws.Range('E30').Value = '  -2.52%  '
$ws.Range('E31').Value = '  -2.52%  '
$ws.Range('D32').Value = '10.29'
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('B33').Value = 'VeChain'
$ws.Range('C33').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D33').Value = '0.0473'
$ws.Range('E33').Value = '  +5.06%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').Value = '35.21'
$ws.Range('E34').Value = '  +3.14%  '
$ws.Range('E35').Value = '  +0.34%  '
$ws.Range('D36').Value = '50.25'
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('D38').Value = '3.35'
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').Value = '0.300'
$ws.Range('E39').Value = '  +7.09%  '
$ws.Range('E40').Value = '  +0.75%  '
$ws.Range('D41').Value = '128.54'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('D43').Value = '16.51'
$ws.Range('E43').Value = '  -4.19%  '
$ws.Range('D44').Value = '2.50'
$ws.Range('E44').Value = '  -3.14%  '
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('D46').Value = '22.25'
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('E47').Value = '  +3.38%  '
$ws.Range('E48').Value = '  +1.60%  '
$ws.Range('D49').Value = '2.051.71'
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('B50').Value = 'BEAM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D50').Value = '0.0326'
$ws.Range('E50').Value = '  +1.82%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.894'
$ws.Range('E51').Value = '  +14.10%  '
